# IGCC Netting Flows Historical - roll the data window forward by two days.
# 17.02.2026 -> 19.02.2026, 18.02.2026 -> 20.02.2026 (both the Timestamp
# serials in column A and the text "Lookup" values in column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 201

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $tsCell = $ws.Cells.Item($r, 1)
    $lookupCell = $ws.Cells.Item($r, 5)

    $ts = $tsCell.Value2
    $tsCell.Value = $ts + 2

    $lookup = $lookupCell.Value2
    $lookup = $lookup.Replace("17.02.2026", "19.02.2026")
    $lookup = $lookup.Replace("18.02.2026", "20.02.2026")
    $lookupCell.Value = $lookup
}
